# Update stats for 2025-12 (row 25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6484
$ws.Range("C25").Value = 1007
$ws.Range("D25").Value = 6043377
$ws.Range("E25").Value = 932.0445712523134
$ws.Range("F25").Value = 10.06620268205738
$ws.Range("G25").Value = 7.356076759061825
$ws.Range("H25").Value = 26.56695105978286
